$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2020" column (Q) that mirrors the existing "2019" column (P):
# copy formats from P3/P4 into Q3/Q4, then set the new values.
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)
$ws.Range("Q3").Value = 2020

$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = [double]"6.7156049127444606E-2"

# The whole data row (now D4:Q4) switches from the custom "0.0" number
# format to the built-in "0.00" format.
$ws.Range("D4:Q4").NumberFormat = "0.00"

# Clear the clipboard marching-ants / leftover selection state and make
# sure the sheet's stored selection is reset to the top-left cell.
$excel.CutCopyMode = 0
$ws.Range("A1").Select()
